# Updates the cryptos list (Price/Volume(1h) columns) for Mon Sep 25 18:29:17 UTC 2023 run.
# Rows 48/49 additionally swap coin identity (name + link) along with their price/volume.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.357.22"
$ws.Range("E2").Value = "  -1.28%  "
# Row 3
$ws.Range("D3").Value = "1.593.12"
$ws.Range("E3").Value = "  -0.23%  "
# Row 4
$ws.Range("E4").Value = "  -0.71%  "
# Row 5
$ws.Range("D5").Value = "'210.13"
$ws.Range("E5").Value = "  -0.82%  "
# Row 6
$ws.Range("D6").Value = "'0.506"
$ws.Range("E6").Value = "  -1.25%  "
# Row 7
$ws.Range("E7").Value = "  -0.76%  "
# Row 8
$ws.Range("E8").Value = "  -0.84%  "
# Row 9
$ws.Range("E9").Value = "  -0.35%  "
# Row 10
$ws.Range("E10").Value = "  -0.51%  "
# Row 11
$ws.Range("D11").Value = "'0.0845"
$ws.Range("E11").Value = "  -0.58%  "
# Row 12
$ws.Range("D12").Value = "1.818.12"
$ws.Range("E12").Value = "  -0.23%  "
# Row 13
$ws.Range("D13").Value = "1.596.70"
$ws.Range("E13").Value = "  -0.10%  "
# Row 14
$ws.Range("E14").Value = "  +0.35%  "
# Row 15
$ws.Range("D15").Value = "'0.517"
$ws.Range("E15").Value = "  -1.54%  "
# Row 16
$ws.Range("D16").Value = "'64.58"
$ws.Range("E16").Value = "  -0.72%  "
# Row 17
$ws.Range("D17").Value = "26.366.88"
$ws.Range("E17").Value = "  -1.16%  "
# Row 18
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  -1.80%  "
# Row 19
$ws.Range("D19").Value = "'7.49"
$ws.Range("E19").Value = "  +4.72%  "
# Row 20
$ws.Range("D20").Value = "'211.10"
$ws.Range("E20").Value = "  +0.81%  "
# Row 21
$ws.Range("E21").Value = "  -0.53%  "
# Row 22
$ws.Range("D22").Value = "'4.27"
$ws.Range("E22").Value = "  -0.66%  "
# Row 23
$ws.Range("E23").Value = "  -2.95%  "
# Row 24
$ws.Range("E24").Value = "  -1.21%  "
# Row 25
$ws.Range("D25").Value = "'145.22"
$ws.Range("E25").Value = "  +0.72%  "
# Row 26
$ws.Range("E26").Value = "  -0.55%  "
# Row 27
$ws.Range("E27").Value = "  -1.06%  "
# Row 28
$ws.Range("E28").Value = "  -0.72%  "
# Row 29
$ws.Range("E29").Value = "  -0.40%  "
# Row 30
$ws.Range("E30").Value = "  -0.99%  "
# Row 31
$ws.Range("E31").Value = "  -0.18%  "
# Row 32
$ws.Range("E32").Value = "  -1.50%  "
# Row 33
$ws.Range("E33").Value = "  +0.05%  "
# Row 34
$ws.Range("D34").Value = "1.303.36"
$ws.Range("E34").Value = "  +1.98%  "
# Row 35
$ws.Range("E35").Value = "  +3.22%  "
# Row 36
$ws.Range("E36").Value = "  -2.08%  "
# Row 37
$ws.Range("E37").Value = "  -0.63%  "
# Row 38
$ws.Range("E38").Value = "  -0.41%  "
# Row 39
$ws.Range("E39").Value = "  -12.85%  "
# Row 40
$ws.Range("E40").Value = "  -1.54%  "
# Row 41
$ws.Range("E41").Value = "  -0.49%  "
# Row 42
$ws.Range("E42").Value = "  +3.13%  "
# Row 43
$ws.Range("D43").Value = "'62.76"
$ws.Range("E43").Value = "  +0.20%  "
# Row 44
$ws.Range("E44").Value = "  -4.59%  "
# Row 45
$ws.Range("E45").Value = "  -1.92%  "
# Row 46
$ws.Range("D46").Value = "1.728.74"
$ws.Range("E46").Value = "  -0.49%  "
# Row 47
$ws.Range("D47").Value = "'88.08"
$ws.Range("E47").Value = "  -2.69%  "
# Row 48
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("E48").Value = "  +5.58%  "
# Row 49
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.50"
$ws.Range("E49").Value = "  -4.19%  "
# Row 50
$ws.Range("D50").Value = "'0.0983"
$ws.Range("E50").Value = "  -4.17%  "
# Row 51
$ws.Range("E51").Value = "  -1.44%  "
